$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Row 1: remove the redundant duplicated "value" header cells (C1:F1) ---
$ws.Range("C1:F1").ClearContents()

# --- Row 8: "Model" -> "production_function" (keeps "Sigmoid" in B8) ---
$ws.Range("A8").Value = "production_function"

# --- Insert a new row 9 for "L_curve" = 1 (everything below shifts down one row) ---
$ws.Range("A9").EntireRow.Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = "0.00E+00"

# --- Remove the old "Deletion" row (now pushed down to row 17) ---
$ws.Range("A17").EntireRow.Delete()

# --- Sheet-view / selection changes ---
# wt_log2_expression loses the "tabSelected" flag (was active before the edit)
$wsWt = $wb.Worksheets.Item("wt_log2_expression")
$wsWt.Activate()

# optimization_parameters becomes the active / selected tab, with a new selection
$ws.Activate()
$ws.Range("C1:G9").Select()
